$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("C1").Value = "funding_source"
$ws.Range("C2").Value = "NC"
$ws.Range("C3").Select()
